# Generate Report for Handoff
# The "b.md" file has completed its handback round-trip: update the
# Overview sheet and the per-locale (zh-cn / de-de) detail sheets for
# the row that tracks b.md so the workbook reflects the newly generated
# handoff report.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet - row 3 is the b.md entry (E=zh-cn status, F=de-de
# status, G=Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-18 22:37:12"

# ---------------------------------------------------------------------
# zh-cn detail sheet - row 3 is the b.md entry
#   C = Status
#   F = Content Duplicate
#   G = Latest Handoff File
#   H = Latest Handoff Datetime
#   P = Error Detail
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces Excel to store this as literal text (matching
# the workbook's existing convention of True/False as strings, not bools).
# Reset the style afterwards so the quote-prefix marker Excel applies
# doesn't leave a spurious cell-style change behind.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-18 22:37:07"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81989117afd2eca6ccaa77c0ebac7f7f34eef237/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cdd3f1571fc2e24a2dbfc9ab66027cd7bb686c0/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666666

# ---------------------------------------------------------------------
# de-de detail sheet - row 3 is the b.md entry
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-18 22:37:12"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81989117afd2eca6ccaa77c0ebac7f7f34eef237/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cdd3f1571fc2e24a2dbfc9ab66027cd7bb686c0/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.16666666666666
